$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "F2r"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 184.0626906666667
$ws.Cells.Item(2, 8).Value = 552.188072
$ws.Cells.Item(2, 9).Value = 0.6510505751503485
$ws.Cells.Item(2, 10).Value = 0.6510505751503486
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 5.914580333333333
$ws.Cells.Item(2, 14).Value = 17.743741
$ws.Cells.Item(2, 15).Value = 0.07170723223214719
$ws.Cells.Item(2, 16).Value = 0.07170723223214717
$ws.Cells.Item(2, 17).Value = 1088.653570317484
$ws.Cells.Item(2, 18).Value = 9797.882132857352
$ws.Cells.Item(2, 19).Value = 0.04668503478717903
$ws.Cells.Item(2, 20).Value = 0.04668503478717903

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "F2r"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 184.0626906666667
$ws.Cells.Item(3, 8).Value = 552.188072
$ws.Cells.Item(3, 9).Value = 0.6510505751503485
$ws.Cells.Item(3, 10).Value = 0.6510505751503486
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 44.940909
$ws.Cells.Item(3, 14).Value = 134.822727
$ws.Cells.Item(3, 15).Value = 0.5448549206821932
$ws.Cells.Item(3, 16).Value = 0.5448549206821931
$ws.Cells.Item(3, 17).Value = 8271.944631545815
$ws.Cells.Item(3, 18).Value = 74447.50168391234
$ws.Cells.Item(3, 19).Value = 0.3547281094836394
$ws.Cells.Item(3, 20).Value = 0.3547281094836394

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "F2r"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 184.0626906666667
$ws.Cells.Item(4, 8).Value = 552.188072
$ws.Cells.Item(4, 9).Value = 0.6510505751503485
$ws.Cells.Item(4, 10).Value = 0.6510505751503486
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 31.626851
$ws.Cells.Item(4, 14).Value = 94.88055300000001
$ws.Cells.Item(4, 15).Value = 0.3834378470856596
$ws.Cells.Item(4, 16).Value = 0.3834378470856596
$ws.Cells.Item(4, 17).Value = 5821.323292373758
$ws.Cells.Item(4, 18).Value = 52391.90963136382
$ws.Cells.Item(4, 19).Value = 0.2496374308795301
$ws.Cells.Item(4, 20).Value = 0.2496374308795301

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "F2r"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 57.4434
$ws.Cells.Item(5, 8).Value = 172.3302
$ws.Cells.Item(5, 9).Value = 0.2031838091312023
$ws.Cells.Item(5, 10).Value = 0.2031838091312023
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 5.914580333333333
$ws.Cells.Item(5, 14).Value = 17.743741
$ws.Cells.Item(5, 15).Value = 0.07170723223214719
$ws.Cells.Item(5, 16).Value = 0.07170723223214717
$ws.Cells.Item(5, 17).Value = 339.7536039198
$ws.Cells.Item(5, 18).Value = 3057.7824352782
$ws.Cells.Item(5, 19).Value = 0.01456974858718339
$ws.Cells.Item(5, 20).Value = 0.01456974858718339

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "F2r"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 57.4434
$ws.Cells.Item(6, 8).Value = 172.3302
$ws.Cells.Item(6, 9).Value = 0.2031838091312023
$ws.Cells.Item(6, 10).Value = 0.2031838091312023
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 44.940909
$ws.Cells.Item(6, 14).Value = 134.822727
$ws.Cells.Item(6, 15).Value = 0.5448549206821932
$ws.Cells.Item(6, 16).Value = 0.5448549206821931
$ws.Cells.Item(6, 17).Value = 2581.5586120506
$ws.Cells.Item(6, 18).Value = 23234.0275084554
$ws.Cells.Item(6, 19).Value = 0.1107056982080871
$ws.Cells.Item(6, 20).Value = 0.1107056982080871

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "F2r"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 57.4434
$ws.Cells.Item(7, 8).Value = 172.3302
$ws.Cells.Item(7, 9).Value = 0.2031838091312023
$ws.Cells.Item(7, 10).Value = 0.2031838091312023
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 31.626851
$ws.Cells.Item(7, 14).Value = 94.88055300000001
$ws.Cells.Item(7, 15).Value = 0.3834378470856596
$ws.Cells.Item(7, 16).Value = 0.3834378470856596
$ws.Cells.Item(7, 17).Value = 1816.7538527334
$ws.Cells.Item(7, 18).Value = 16350.7846746006
$ws.Cells.Item(7, 19).Value = 0.07790836233593181
$ws.Cells.Item(7, 20).Value = 0.07790836233593179

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "F2r"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 41.21033366666666
$ws.Cells.Item(8, 8).Value = 123.631001
$ws.Cells.Item(8, 9).Value = 0.1457656157184491
$ws.Cells.Item(8, 10).Value = 0.1457656157184491
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 5.914580333333333
$ws.Cells.Item(8, 14).Value = 17.743741
$ws.Cells.Item(8, 15).Value = 0.07170723223214719
$ws.Cells.Item(8, 16).Value = 0.07170723223214717
$ws.Cells.Item(8, 17).Value = 243.7418290349712
$ws.Cells.Item(8, 18).Value = 2193.676461314741
$ws.Cells.Item(8, 19).Value = 0.01045244885778476
$ws.Cells.Item(8, 20).Value = 0.01045244885778475

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "F2r"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 41.21033366666666
$ws.Cells.Item(9, 8).Value = 123.631001
$ws.Cells.Item(9, 9).Value = 0.1457656157184491
$ws.Cells.Item(9, 10).Value = 0.1457656157184491
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 44.940909
$ws.Cells.Item(9, 14).Value = 134.822727
$ws.Cells.Item(9, 15).Value = 0.5448549206821932
$ws.Cells.Item(9, 16).Value = 0.5448549206821931
$ws.Cells.Item(9, 17).Value = 1852.029855173303
$ws.Cells.Item(9, 18).Value = 16668.26869655972
$ws.Cells.Item(9, 19).Value = 0.07942111299046666
$ws.Cells.Item(9, 20).Value = 0.07942111299046664

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Gnai2"
$ws.Cells.Item(10, 3).Value = "F2r"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 41.21033366666666
$ws.Cells.Item(10, 8).Value = 123.631001
$ws.Cells.Item(10, 9).Value = 0.1457656157184491
$ws.Cells.Item(10, 10).Value = 0.1457656157184491
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 31.626851
$ws.Cells.Item(10, 14).Value = 94.88055300000001
$ws.Cells.Item(10, 15).Value = 0.3834378470856596
$ws.Cells.Item(10, 16).Value = 0.3834378470856596
$ws.Cells.Item(10, 17).Value = 1303.35308253595
$ws.Cells.Item(10, 18).Value = 11730.17774282355
$ws.Cells.Item(10, 19).Value = 0.05589205387019772
$ws.Cells.Item(10, 20).Value = 0.05589205387019771
